$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# NumberFormat/Style juggling keeps the cells as plain text (matching the
# original inline-string cells) instead of letting Excel auto-coerce
# numeric-looking strings (e.g. "227.61") into real numbers.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '37.755.18'
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.047.59'
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  +1.06%  '
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '227.61'
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +0.32%  '
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  +0.69%  '
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '59.98'
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  +1.20%  '
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.376'
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -1.71%  '
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0836'
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  +2.98%  '
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +0.80%  '
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '2.350.27'
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  +1.03%  '
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -0.83%  '
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '21.39'
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  +1.78%  '
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +6.62%  '
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  +1.09%  '
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.041.80'
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  +0.55%  '
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '37.726.80'
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '69.41'
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -0.54%  '
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '5.92'
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -1.55%  '
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '222.09'
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -1.08%  '
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  +0.53%  '
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +1.89%  '
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +3.43%  '
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '168.89'
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +2.34%  '
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.32'
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +0.88%  '
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +0.47%  '
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '18.76'
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -0.71%  '
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  +0.58%  '
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  +5.86%  '
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -0.82%  '
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.52'
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +0.96%  '
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  +0.42%  '
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '6.49'
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +2.32%  '
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  +4.41%  '
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.48'
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  +7.84%  '
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '18.26'
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +8.93%  '
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.532.07'
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +0.18%  '
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '97.72'
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  +1.02%  '
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '4.16'
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -0.87%  '
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0892'
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -2.45%  '
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +0.45%  '
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +0.36%  '
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -0.20%  '
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '7.09'
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  +0.00%  '
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.239.38'
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +1.08%  '
$c.Style = "Normal"
